# Auto-generated Excel COM-interop edit script
# Updates cached market-price columns (H-N) on several sheets to match
# the scheduled runner's refreshed values from Ifrit_Profits.xlsx.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 971.93024
$ws.Range("I40").Value = 940.7273
$ws.Range("J40").Value = 1074.9
$ws.Range("K40").Value = 940.7273
$ws.Range("L40").Value = 1074.9
$ws.Range("M40").Value = -765.7273
$ws.Range("N40").Value = -1424.9
$ws.Range("H132").Value = 316391.03
$ws.Range("I132").Value = 337403.78
$ws.Range("K132").Value = 1012211.34
$ws.Range("M132").Value = -1009681.34
$ws.Range("H135").Value = 3286.3076
$ws.Range("I135").Value = 1281.5555
$ws.Range("J135").Value = 7797
$ws.Range("K135").Value = 11533.9995
$ws.Range("L135").Value = 70173
$ws.Range("M135").Value = -8998.9995
$ws.Range("N135").Value = -75243
$ws.Range("H137").Value = 2619.4375
$ws.Range("I137").Value = 1375.9166
$ws.Range("J137").Value = 6350
$ws.Range("K137").Value = 4127.7498
$ws.Range("L137").Value = 19050
$ws.Range("M137").Value = -1577.7498
$ws.Range("N137").Value = -24150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4114.8335
$ws.Range("I61").Value = 2535.7144
$ws.Range("J61").Value = 6325.6
$ws.Range("K61").Value = 2535.7144
$ws.Range("L61").Value = 6325.6
$ws.Range("M61").Value = -2323.7144
$ws.Range("N61").Value = -6749.6
$ws.Range("H97").Value = 562.2222
$ws.Range("I97").Value = 510
$ws.Range("J97").Value = 666.6667
$ws.Range("K97").Value = 510
$ws.Range("L97").Value = 666.6667
$ws.Range("M97").Value = -14
$ws.Range("N97").Value = -1658.6667
$ws.Range("H132").Value = 3726.6316
$ws.Range("I132").Value = 3557.7144
$ws.Range("K132").Value = 10673.1432
$ws.Range("M132").Value = -8143.143199999999
$ws.Range("H136").Value = 4114.8335
$ws.Range("I136").Value = 2535.7144
$ws.Range("J136").Value = 6325.6
$ws.Range("K136").Value = 7607.1432
$ws.Range("L136").Value = 18976.8
$ws.Range("M136").Value = -5057.1432
$ws.Range("N136").Value = -24076.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1507.9615
$ws.Range("I20").Value = 1311.4
$ws.Range("J20").Value = 1776
$ws.Range("K20").Value = 1311.4
$ws.Range("L20").Value = 1776
$ws.Range("M20").Value = -1064.4
$ws.Range("N20").Value = -2270
$ws.Range("H134").Value = 2250
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 6750
$ws.Range("N134").Value = -11820
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1684.421
$ws.Range("I31").Value = 1000.3077
$ws.Range("J31").Value = 3166.6667
$ws.Range("K31").Value = 1000.3077
$ws.Range("L31").Value = 3166.6667
$ws.Range("M31").Value = -705.3077
$ws.Range("N31").Value = -3756.6667
$ws.Range("H34").Value = 1684.421
$ws.Range("I34").Value = 1000.3077
$ws.Range("J34").Value = 3166.6667
$ws.Range("K34").Value = 1000.3077
$ws.Range("L34").Value = 3166.6667
$ws.Range("M34").Value = -798.3077
$ws.Range("N34").Value = -3570.6667
$ws.Range("H58").Value = 2567.087
$ws.Range("I58").Value = 1918.3846
$ws.Range("J58").Value = 3410.4
$ws.Range("K58").Value = 1918.3846
$ws.Range("L58").Value = 3410.4
$ws.Range("M58").Value = -1715.3846
$ws.Range("N58").Value = -3816.4
$ws.Range("H105").Value = 1027.3
$ws.Range("I105").Value = 1009.5714
$ws.Range("J105").Value = 1068.6666
$ws.Range("K105").Value = 1009.5714
$ws.Range("L105").Value = 1068.6666
$ws.Range("M105").Value = 737.4286
$ws.Range("N105").Value = -4562.6666
$ws.Range("H133").Value = 26969.445
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 26969.445
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 26969.445
$ws.Range("N133").Value = -32029.445
$ws.Range("H134").Value = 2632.8572
$ws.Range("I134").Value = 2641.4814
$ws.Range("K134").Value = 7924.4442
$ws.Range("M134").Value = -5389.4442
$ws.Range("H136").Value = 2567.087
$ws.Range("I136").Value = 1918.3846
$ws.Range("J136").Value = 3410.4
$ws.Range("K136").Value = 5755.1538
$ws.Range("L136").Value = 10231.2
$ws.Range("M136").Value = -3205.1538
$ws.Range("N136").Value = -15331.2
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1880.091
$ws.Range("I3").Value = 931.5
$ws.Range("K3").Value = 2794.5
$ws.Range("M3").Value = -2682.5
$ws.Range("H117").Value = 1753.0714
$ws.Range("I117").Value = 811.4286
$ws.Range("J117").Value = 2694.7144
$ws.Range("K117").Value = 2434.2858
$ws.Range("L117").Value = 8084.1432
$ws.Range("M117").Value = 1007.7142
$ws.Range("N117").Value = -14968.1432
$ws.Range("H129").Value = 1125.2307
$ws.Range("I129").Value = 659.6667
$ws.Range("J129").Value = 1264.9
$ws.Range("K129").Value = 1979.0001
$ws.Range("L129").Value = 3794.7
$ws.Range("M129").Value = 3020.9999
$ws.Range("N129").Value = -13794.7
$ws.Range("H131").Value = 1641973.1
$ws.Range("I131").Value = 10296.667
$ws.Range("J131").Value = 1819974.1
$ws.Range("K131").Value = 30890.001
$ws.Range("L131").Value = 5459922.300000001
$ws.Range("M131").Value = -25850.001
$ws.Range("N131").Value = -5470002.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H126").Value = 1362.625
$ws.Range("I126").Value = 1483.5
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 4450.5
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -1980.5
$ws.Range("N126").Value = -7940
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 20233.334
$ws.Range("I136").Value = 34466.668
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 103400.004
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -100850.004
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17958.666
$ws.Range("I136").Value = 29743.428
$ws.Range("J136").Value = 1460
$ws.Range("K136").Value = 89230.284
$ws.Range("L136").Value = 4380
$ws.Range("M136").Value = -86680.284
$ws.Range("N136").Value = -9480
